$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.525.33'
$ws.Range('E2').Value = '  +2.84%  '

$ws.Range('D3').Value = '2.076.38'
$ws.Range('E3').Value = '  +3.80%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.05%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.618'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.33%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.37'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +6.67%  '

$ws.Range('E8').Value = '  +0.02%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.385'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.04%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '59.16'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.54%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0763'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.18%  '

$ws.Range('E12').Value = '  +3.85%  '

$ws.Range('D13').Value = '2.377.11'
$ws.Range('E13').Value = '  +3.70%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.59'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.88%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.15'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.20%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.781'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.28%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.20'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.74%  '

$ws.Range('D18').Value = '2.071.90'
$ws.Range('E18').Value = '  +3.48%  '

$ws.Range('D19').Value = '37.700.75'
$ws.Range('E19').Value = '  +3.49%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.24'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +18.26%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '70.18'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.54%  '

$ws.Range('D22').Value = '0.0₃0817'
$ws.Range('E22').Value = '  +1.58%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '227.08'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.50%  '

$ws.Range('E24').Value = '  +0.04%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.45'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.48%  '

$ws.Range('E26').Value = '  +1.17%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.72'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.38%  '

$ws.Range('E28').Value = '  +9.92%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.94'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.20%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.32'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.35%  '

$ws.Range('E31').Value = '  +1.63%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.119'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.52%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.52'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.42%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0625'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.39%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.58'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +8.34%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.57'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.36%  '

$ws.Range('B37').Value = 'BinanceUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.07%  '

$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.36'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.50%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.79'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.49%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.87'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.54%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.61'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +26.17%  '

$ws.Range('E42').Value = '  -1.27%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0954'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.10%  '

$ws.Range('E44').Value = '  +7.97%  '

$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '96.01'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.35%  '

$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '1.457.98'
$ws.Range('E46').Value = '  +0.42%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0213'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.16%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '15.85'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.57%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.04'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.51%  '

$ws.Range('E50').Value = '  +6.22%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.93'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.73%  '
